# Add a new student record (row 4) to the roster, matching the same
# layout / hyperlink-on-email pattern used by the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Sidhartha Nambiar"
$ws.Range("B4").Value = "1MS16CS044"
$ws.Range("C4").Value = 8848779798
$ws.Range("D4").Value = "nambiar.sidhartha00@gmail.com"
$ws.Range("E4").Value = "CSE"
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = "Sidhrtha Appa"
$ws.Range("H4").Value = "Sidhartha Amma"

# Mirror the mailto hyperlinks already present on D2/D3.
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:nambiar.sidhartha00@gmail.com")
$ws.Range("D4").Style = "Hyperlink"

# Column A needs to widen to fit the new, longer name.
[void]$ws.Columns.Item(1).AutoFit()

# Leave the selection on the last-entered cell, as happens after typing.
[void]$ws.Range("H4").Select()
